$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# B6 previously held "proceso" (shared string) - now it should show "Leido"
$ws.Range("B6").Value = "Leido"

# B7 is a new cell holding the renamed string (previously "proceso", now "Preceso")
$ws.Range("B7").Value = "Preceso"

# Update the active selection to B7
$ws.Range("B7").Select()
